{"js": "// Remove the \"Attn: Designated Copyright Agent\" line from the notice\n// address block (the paragraph that sits between \"Roil Technology, LLC\"\n// and \"Jonathan Levy\").\nconst body = context.document.body;\n\nconst results = body.search(\"Attn: Designated Copyright Agent\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  const para = results.items[i].paragraphs.getFirst();\n  para.delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the \"Attn: Designated Copyright Agent\" line from the notice\n# address block (the paragraph that sits between \"Roil Technology, LLC\"\n# and \"Jonathan Levy\").\n$d = $word.ActiveDocument\n$range = $d.Content\n\n$range.Find.ClearFormatting()\n$range.Find.Text = \"Attn: Designated Copyright Agent\"\n$range.Find.MatchCase = $true\n$range.Find.MatchWholeWord = $false\n$range.Find.Wrap = 1  # wdFindContinue\n\nwhile ($range.Find.Execute()) {\n    # Expand the found text to the whole paragraph (including its mark)\n    # and remove it entirely.\n    $range.Expand(4) | Out-Null   # wdParagraph\n    $range.Delete()\n    $range.Collapse(0)            # wdCollapseEnd\n}\n"}
